$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# Set the new header cell value
$ws.Range("A1").Value = "kolom"

# Move the active selection to A2, matching the post-edit workbook state
$ws.Range("A2").Select()
